$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Repulled dSF (F column) values per row
$ws.Range("F2").Value = -5
$ws.Range("F3").Value = -1
$ws.Range("F5").Value = 1
$ws.Range("F6").Value = -1
$ws.Range("F7").Value = -4
$ws.Range("F8").Value = -4
$ws.Range("F10").Value = -7
$ws.Range("F15").Value = -2
$ws.Range("F19").Value = 5
$ws.Range("F28").Value = -1
$ws.Range("F29").Value = -1
$ws.Range("F30").Value = -7
$ws.Range("F31").Value = 2
$ws.Range("F32").Value = 1
$ws.Range("F33").Value = 5
$ws.Range("F35").Value = 1
$ws.Range("F36").Value = 2
$ws.Range("F37").Value = 1
$ws.Range("F38").Value = -2
$ws.Range("F40").Value = -1
$ws.Range("F41").Value = 1
$ws.Range("F42").Value = -2
